$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 650.375
$ws.Range("I15").Value = 650.375
$ws.Range("K15").Value = 1951.125
$ws.Range("M15").Value = -1782.125
$ws.Range("H19").Value = 716.36365
$ws.Range("I19").Value = 150
$ws.Range("K19").Value = 150
$ws.Range("M19").Value = 25
$ws.Range("H40").Value = 4032.9285
$ws.Range("I40").Value = 4026.4
$ws.Range("K40").Value = 4026.4
$ws.Range("M40").Value = -3851.4
$ws.Range("H43").Value = 4856.1113
$ws.Range("I43").Value = 1801
$ws.Range("K43").Value = 1801
$ws.Range("M43").Value = -1732
$ws.Range("H96").Value = 239.66667
$ws.Range("I96").Value = 262.54544
$ws.Range("J96").Value = 176.75
$ws.Range("K96").Value = 787.63632
$ws.Range("L96").Value = 530.25
$ws.Range("M96").Value = 585.36368
$ws.Range("N96").Value = -3276.25
$ws.Range("H116").Value = 5200.2188
$ws.Range("I116").Value = 4103.222
$ws.Range("K116").Value = 4103.222
$ws.Range("M116").Value = -661.2219999999998
$ws.Range("H129").Value = 125001580
$ws.Range("I129").Value = 166667780
$ws.Range("K129").Value = 500003340
$ws.Range("M129").Value = -499998340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10002.129
$ws.Range("I32").Value = 5636.978
$ws.Range("K32").Value = 5636.978
$ws.Range("M32").Value = -5349.978
$ws.Range("H47").Value = 18000
$ws.Range("J47").Value = 18000
$ws.Range("L47").Value = 18000
$ws.Range("N47").Value = -19450
$ws.Range("H97").Value = 1546278
$ws.Range("I97").Value = 2941969.2
$ws.Range("K97").Value = 2941969.2
$ws.Range("M97").Value = -2941473.2
$ws.Range("H102").Value = 5560321.5
$ws.Range("I102").Value = 8337679.5
$ws.Range("K102").Value = 8337679.5
$ws.Range("M102").Value = -8336057.5
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 2569.907
$ws.Range("I132").Value = 2107.5134
$ws.Range("K132").Value = 6322.540199999999
$ws.Range("M132").Value = -3792.540199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 833.5
$ws.Range("J11").Value = 500.2
$ws.Range("L11").Value = 500.2
$ws.Range("N11").Value = -780.2
$ws.Range("H86").Value = 16668532
$ws.Range("I86").Value = 33334900
$ws.Range("J86").Value = 2163.3333
$ws.Range("K86").Value = 33334900
$ws.Range("L86").Value = 2163.3333
$ws.Range("M86").Value = -33333777
$ws.Range("N86").Value = -4409.3333
$ws.Range("H89").Value = 16668532
$ws.Range("I89").Value = 33334900
$ws.Range("J89").Value = 2163.3333
$ws.Range("K89").Value = 166674500
$ws.Range("L89").Value = 10816.6665
$ws.Range("M89").Value = -166668884
$ws.Range("N89").Value = -22048.6665
$ws.Range("H94").Value = 3455706.2
$ws.Range("I94").Value = 5883060
$ws.Range("K94").Value = 5883060
$ws.Range("M94").Value = -5882609
$ws.Range("H99").Value = 6214913.5
$ws.Range("I99").Value = 10207651
$ws.Range("K99").Value = 10207651
$ws.Range("M99").Value = -10206153
$ws.Range("H107").Value = 2234745
$ws.Range("I107").Value = 3107515
$ws.Range("K107").Value = 3107515
$ws.Range("M107").Value = -3105595
$ws.Range("H134").Value = 3234.4656
$ws.Range("I134").Value = 1147.75
$ws.Range("J134").Value = 7871.6113
$ws.Range("K134").Value = 3443.25
$ws.Range("L134").Value = 23614.8339
$ws.Range("M134").Value = -908.25
$ws.Range("N134").Value = -28684.8339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1614.6154
$ws.Range("I16").Value = 1363.7273
$ws.Range("K16").Value = 1363.7273
$ws.Range("M16").Value = -1076.7273
$ws.Range("H22").Value = 727.2308
$ws.Range("J22").Value = 1038.25
$ws.Range("L22").Value = 1038.25
$ws.Range("N22").Value = -1738.25
$ws.Range("H29").Value = 26010.5
$ws.Range("J29").Value = 26010.5
$ws.Range("L29").Value = 26010.5
$ws.Range("N29").Value = -26596.5
$ws.Range("H31").Value = 22962.84
$ws.Range("I31").Value = 2927
$ws.Range("K31").Value = 2927
$ws.Range("M31").Value = -2632
$ws.Range("H34").Value = 22962.84
$ws.Range("I34").Value = 2927
$ws.Range("K34").Value = 2927
$ws.Range("M34").Value = -2725
$ws.Range("H58").Value = 4249.2905
$ws.Range("I58").Value = 5269.875
$ws.Range("J58").Value = 3160.6667
$ws.Range("K58").Value = 5269.875
$ws.Range("L58").Value = 3160.6667
$ws.Range("M58").Value = -5066.875
$ws.Range("N58").Value = -3566.6667
$ws.Range("H99").Value = 3220.5557
$ws.Range("I99").Value = 2726.7856
$ws.Range("J99").Value = 4948.75
$ws.Range("K99").Value = 2726.7856
$ws.Range("L99").Value = 4948.75
$ws.Range("M99").Value = -1228.7856
$ws.Range("N99").Value = -7944.75
$ws.Range("H105").Value = 1109.5
$ws.Range("I105").Value = 1109.5
$ws.Range("K105").Value = 1109.5
$ws.Range("M105").Value = 637.5
$ws.Range("H113").Value = 1614.6154
$ws.Range("I113").Value = 1363.7273
$ws.Range("K113").Value = 1363.7273
$ws.Range("M113").Value = 806.2727
$ws.Range("H117").Value = 64997.5
$ws.Range("J117").Value = 64997.5
$ws.Range("L117").Value = 64997.5
$ws.Range("N117").Value = -74175.5
$ws.Range("H126").Value = 3220.5557
$ws.Range("I126").Value = 2726.7856
$ws.Range("J126").Value = 4948.75
$ws.Range("K126").Value = 8180.3568
$ws.Range("L126").Value = 14846.25
$ws.Range("M126").Value = -5710.3568
$ws.Range("N126").Value = -19786.25
$ws.Range("H132").Value = 47268.535
$ws.Range("I132").Value = 31390.428
$ws.Range("J132").Value = 116735.25
$ws.Range("K132").Value = 94171.284
$ws.Range("L132").Value = 350205.75
$ws.Range("M132").Value = -91641.284
$ws.Range("N132").Value = -355265.75
$ws.Range("H134").Value = 2982.878
$ws.Range("I134").Value = 2045.2258
$ws.Range("K134").Value = 6135.6774
$ws.Range("M134").Value = -3600.6774
$ws.Range("H136").Value = 4249.2905
$ws.Range("I136").Value = 5269.875
$ws.Range("J136").Value = 3160.6667
$ws.Range("K136").Value = 15809.625
$ws.Range("L136").Value = 9482.000100000001
$ws.Range("M136").Value = -13259.625
$ws.Range("N136").Value = -14582.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2165603.5
$ws.Range("J97").Value = 1389
$ws.Range("L97").Value = 1389
$ws.Range("N97").Value = -2381
$ws.Range("H108").Value = 36929
$ws.Range("I108").Value = 39375.5
$ws.Range("J108").Value = 33667
$ws.Range("K108").Value = 39375.5
$ws.Range("L108").Value = 33667
$ws.Range("M108").Value = -35535.5
$ws.Range("N108").Value = -41347
$ws.Range("H126").Value = 3581618.2
$ws.Range("I126").Value = 5053143.5
$ws.Range("J126").Value = 3091110
$ws.Range("K126").Value = 15159430.5
$ws.Range("L126").Value = 9273330
$ws.Range("M126").Value = -15156960.5
$ws.Range("N126").Value = -9278270
$ws.Range("H132").Value = 3450.7646
$ws.Range("I132").Value = 3226.1072
$ws.Range("K132").Value = 9678.321599999999
$ws.Range("M132").Value = -7148.321599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 38240
$ws.Range("J6").Value = 38240
$ws.Range("L6").Value = 38240
$ws.Range("N6").Value = -38464
$ws.Range("H16").Value = 779.6
$ws.Range("J16").Value = 848.5
$ws.Range("L16").Value = 848.5
$ws.Range("N16").Value = -1188.5
$ws.Range("H88").Value = 18999.666
$ws.Range("I88").Value = 11000
$ws.Range("J88").Value = 22999.5
$ws.Range("K88").Value = 11000
$ws.Range("L88").Value = 22999.5
$ws.Range("N88").Value = -23855.5
$ws.Range("M88").Value = -10572
$ws.Range("H91").Value = 18999.666
$ws.Range("I91").Value = 11000
$ws.Range("J91").Value = 22999.5
$ws.Range("K91").Value = 11000
$ws.Range("L91").Value = 22999.5
$ws.Range("N91").Value = -25963.5
$ws.Range("M91").Value = -9518
$ws.Range("H93").Value = 37051430
$ws.Range("I93").Value = 55557816
$ws.Range("K93").Value = 55557816
$ws.Range("M93").Value = -55556568
$ws.Range("H95").Value = 24499.5
$ws.Range("J95").Value = 24499.5
$ws.Range("L95").Value = 24499.5
$ws.Range("N95").Value = -29991.5
$ws.Range("H103").Value = 23933.166
$ws.Range("J103").Value = 23933.166
$ws.Range("L103").Value = 23933.166
$ws.Range("N103").Value = -26277.166
$ws.Range("H138").Value = 82123.5
$ws.Range("J138").Value = 82123.5
$ws.Range("L138").Value = 82123.5
$ws.Range("N138").Value = -92403.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 23000
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 23000
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H93").Value = 49999
$ws.Range("J93").Value = 49999
$ws.Range("L93").Value = 49999
$ws.Range("N93").Value = -54991
$ws.Range("H132").Value = 16147809
$ws.Range("I132").Value = 18523242
$ws.Range("J132").Value = 113640.875
$ws.Range("K132").Value = 55569726
$ws.Range("L132").Value = 340922.625
$ws.Range("M132").Value = -55567196
$ws.Range("N132").Value = -345982.625
